$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number -> @(Coin, Link, Price, Volume(1h))
# Only rows whose B/C/D/E values actually change vs. the original are listed here.
$rows = @{
    2  = @($null, $null, "65.082.69", "  -2.37%  ")
    3  = @($null, $null, "3.468.79", "  -1.69%  ")
    4  = @($null, $null, $null, "  +0.08%  ")
    5  = @($null, $null, "586.49", "  -3.41%  ")
    6  = @($null, $null, "136.58", "  -4.72%  ")
    7  = @($null, $null, "3.467.84", "  -1.69%  ")
    8  = @($null, $null, $null, "  +0.11%  ")
    9  = @($null, $null, "0.488", "  -4.25%  ")
    10 = @($null, $null, $null, "  -6.53%  ")
    11 = @($null, $null, "7.15", "  -7.23%  ")
    12 = @($null, $null, "0.379", "  -6.59%  ")
    13 = @($null, $null, "4.061.89", "  -1.41%  ")
    14 = @($null, $null, "0.0000181", "  -7.26%  ")
    15 = @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.448.88", "  -2.16%  ")
    16 = @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "26.37", "  -8.14%  ")
    18 = @($null, $null, "65.034.05", "  -2.14%  ")
    19 = @($null, $null, "9.63", "  -10.86%  ")
    20 = @($null, $null, "5.73", "  -6.46%  ")
    21 = @($null, $null, "13.85", "  -5.45%  ")
    22 = @($null, $null, "388.52", "  -8.28%  ")
    23 = @($null, $null, "0.552", "  -6.35%  ")
    24 = @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.08%  ")
    25 = @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "72.45", "  -5.91%  ")
    26 = @("LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "5.76", "  -0.14%  ")
    27 = @($null, $null, "3.613.72", "  -1.64%  ")
    28 = @($null, $null, "0.0000109", "  -4.66%  ")
    29 = @($null, $null, "1.00", "  +0.32%  ")
    30 = @($null, $null, "7.33", "  -7.46%  ")
    31 = @($null, $null, "8.16", "  -8.82%  ")
    32 = @($null, $null, "2.21", "  -10.27%  ")
    33 = @($null, $null, "3.488.16", "  -1.30%  ")
    35 = @($null, $null, $null, "  -6.96%  ")
    36 = @($null, $null, "22.99", "  -5.17%  ")
    37 = @($null, $null, "170.38", "  -1.73%  ")
    38 = @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.82", "  -9.94%  ")
    39 = @("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.19", "  -10.48%  ")
    40 = @($null, $null, "1.45", "  -11.05%  ")
    41 = @($null, $null, "4.72", "  -9.41%  ")
    42 = @($null, $null, "0.0773", "  -4.78%  ")
    43 = @($null, $null, "0.810", "  -5.02%  ")
    44 = @($null, $null, $null, "  +0.22%  ")
    45 = @($null, $null, "42.37", "  -6.90%  ")
    46 = @($null, $null, "24.65", "  +7.20%  ")
    47 = @($null, $null, "4.33", "  -13.34%  ")
    48 = @($null, $null, $null, "  -9.94%  ")
    49 = @($null, $null, $null, "  +1.49%  ")
    50 = @($null, $null, "6.64", "  -6.12%  ")
    51 = @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.225.18", "  -3.84%  ")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    if ($null -ne $vals[0]) { $ws.Cells.Item($r, 2).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($r, 3).Value = $vals[1] }
    if ($null -ne $vals[2]) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[2]
    }
    if ($null -ne $vals[3]) { $ws.Cells.Item($r, 5).Value = $vals[3] }
}
